$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing data row (row 3, columns A:F) into new rows 7-15,
# matching the report's new larger result set.
$source = $ws.Range("A3:F3")
$source.Copy()

for ($r = 7; $r -le 15; $r++) {
    $target = $ws.Range("A" + $r + ":F" + $r)
    $target.PasteSpecial()
}
